$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph near the end of the document.
$count = $d.Paragraphs.Count
$target = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "Requisitos") {
        $target = $i
    }
}

if ($target -ge 1) {
    # Remove the "Requisitos" heading paragraph and everything after it
    # (the "LOM3074 -  Processamento de Cerâmicas II  (Requisito)" bullet).
    $startPara = $d.Paragraphs($target)
    $endPara = $d.Paragraphs($d.Paragraphs.Count)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
